# Update row 9 (Ano 2025) figures in faturamento_anual sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 4236794.19
$ws.Range("C9").Value = 667716.24
$ws.Range("D9").Value = 4904510.430000001
$ws.Range("E9").Value = 13.61433010552288
$ws.Range("F9").Value = 86.38566989447712
$ws.Range("G9").Value = -35.46855479066289
$ws.Range("H9").Value = -23.48935319336804
$ws.Range("I9").Value = 42605
$ws.Range("J9").Value = 1835
$ws.Range("K9").Value = 44440
$ws.Range("L9").Value = 30876
$ws.Range("M9").Value = 158.8453954527789
$ws.Range("N9").Value = 8.44698679564473
